$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

# Copy the format from A2 (Hyperlink-derived style) onto B2, matching the
# style Excel used when this row was filled in
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in row 2 (URL / UserName / Password) to match the header columns in row 1
$ws.Range("A2").Value = "https://emgi-dev4.login.ca3.oraclecloud.com"
$ws.Range("B2").Value = "Anusuya.Lakkannai@snb.ca"
$ws.Range("C2").Value = "Aishu@123"

# Update selection to B9 as recorded in the saved file
$ws.Range("B9").Select()
